$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Quyết định thôi học" (column U) values for rows 2-11
# into the "Quyết định khen thưởng" (column W) values, clearing the
# original U cells.
for ($row = 2; $row -le 11; $row++) {
    $srcCell = $ws.Range("U$row")
    $dstCell = $ws.Range("W$row")
    $value = $srcCell.Value2
    if ($null -ne $value -and $value -ne "") {
        $dstCell.Value = $value
        $srcCell.ClearContents()
    }
}
